# Auto-generated: update cryptos Price (D) and Volume(1h) (E) columns
# to match the values captured in the Sat Aug 10 23:13:23 UTC 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.789.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "'2.597.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'523.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.95%  "
$ws.Range("D6").Value = "'154.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").Value = "'6.67"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").Value = "'0.346"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "'3.053.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "'60.807.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").Value = "'21.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "'2.602.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "'354.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.47%  "
$ws.Range("D20").Value = "'10.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'60.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").Value = "'0.425"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "'2.717.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "'0.0₃0842"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "'6.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.20%  "
$ws.Range("D32").Value = "'19.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("D34").Value = "'149.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("D35").Value = "'4.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.48%  "
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("D37").Value = "'0.915"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.73%  "
$ws.Range("D38").Value = "'0.908"
$ws.Range("D38").Style = "Normal"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").Value = "'3.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").Value = "'36.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("D42").Value = "'291.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("D44").Value = "'0.624"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").Value = "'0.0558"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "'19.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("D48").Value = "'4.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("D50").Value = "'10.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").Value = "'19.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.80%  "
